# Auto-generated edit script applying the Hyperion_Profits market-data refresh
# (regenerated pricing figures in columns H:N across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 236.16667
$ws.Range("I5").Value = 279.6
$ws.Range("J5").Value = 19
$ws.Range("K5").Value = 279.6
$ws.Range("L5").Value = 19
$ws.Range("M5").Value = -164.6
$ws.Range("N5").Value = -249

$ws.Range("H13").Value = 594.25
$ws.Range("I13").Value = 405
$ws.Range("J13").Value = 657.3333
$ws.Range("K13").Value = 405
$ws.Range("L13").Value = 657.3333
$ws.Range("N13").Value = -995.3333
$ws.Range("M13").Value = -236

$ws.Range("H62").Value = 5016.72
$ws.Range("I62").Value = 3714.2856
$ws.Range("J62").Value = 6674.364
$ws.Range("K62").Value = 3714.2856
$ws.Range("L62").Value = 6674.364
$ws.Range("M62").Value = -3090.2856
$ws.Range("N62").Value = -7922.364

$ws.Range("H65").Value = 5016.72
$ws.Range("I65").Value = 3714.2856
$ws.Range("J65").Value = 6674.364
$ws.Range("K65").Value = 18571.428
$ws.Range("L65").Value = 33371.82
$ws.Range("M65").Value = -15451.428
$ws.Range("N65").Value = -39611.82

$ws.Range("H106").Value = 62501610
$ws.Range("I106").Value = 71430020
$ws.Range("K106").Value = 71430020
$ws.Range("M106").Value = -71429389

$ws.Range("H113").Value = 8297.5
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 8297.5
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 8297.5
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -14805.5

$ws.Range("H116").Value = 5695.857
$ws.Range("I116").Value = 4996
$ws.Range("J116").Value = 6220.75
$ws.Range("K116").Value = 4996
$ws.Range("L116").Value = 6220.75
$ws.Range("M116").Value = -1554
$ws.Range("N116").Value = -13104.75

$ws.Range("H136").Value = 60000
$ws.Range("J136").Value = 60000
$ws.Range("L136").Value = 60000
$ws.Range("N136").Value = -70200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2927.6
$ws.Range("I32").Value = 2163.761
$ws.Range("J32").Value = 11711.75
$ws.Range("K32").Value = 2163.761
$ws.Range("L32").Value = 11711.75
$ws.Range("M32").Value = -1876.761
$ws.Range("N32").Value = -12285.75

$ws.Range("H74").Value = 54424.35
$ws.Range("I74").Value = 30823.871
$ws.Range("J74").Value = 298296
$ws.Range("K74").Value = 30823.871
$ws.Range("L74").Value = 298296
$ws.Range("M74").Value = -29949.871
$ws.Range("N74").Value = -300044

$ws.Range("H77").Value = 54424.35
$ws.Range("I77").Value = 30823.871
$ws.Range("J77").Value = 298296
$ws.Range("K77").Value = 154119.355
$ws.Range("L77").Value = 1491480
$ws.Range("M77").Value = -149751.355
$ws.Range("N77").Value = -1500216

$ws.Range("H122").Value = 1897567.2
$ws.Range("I122").Value = 3134.25
$ws.Range("J122").Value = 2980100.5
$ws.Range("K122").Value = 9402.75
$ws.Range("L122").Value = 8940301.5
$ws.Range("M122").Value = -6952.75
$ws.Range("N122").Value = -8945201.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 838.8
$ws.Range("I22").Value = 838.8
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 838.8
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -665.8
$ws.Range("N22").ClearContents()

$ws.Range("H140").Value = 79694
$ws.Range("J140").Value = 79694
$ws.Range("L140").Value = 79694
$ws.Range("N140").Value = -90054

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 266.05884
$ws.Range("I7").Value = 202.94118
$ws.Range("J7").Value = 329.17648
$ws.Range("K7").Value = 202.94118
$ws.Range("L7").Value = 329.17648
$ws.Range("M7").Value = -89.94118
$ws.Range("N7").Value = -555.1764800000001

$ws.Range("H31").Value = 25134.605
$ws.Range("I31").Value = 3938.15
$ws.Range("K31").Value = 3938.15
$ws.Range("M31").Value = -3643.15

$ws.Range("H34").Value = 25134.605
$ws.Range("I34").Value = 3938.15
$ws.Range("K34").Value = 3938.15
$ws.Range("M34").Value = -3736.15

$ws.Range("H58").Value = 1412.625
$ws.Range("I58").Value = 786.4138
$ws.Range("K58").Value = 786.4138
$ws.Range("M58").Value = -583.4138

$ws.Range("H99").Value = 4511.933
$ws.Range("I99").Value = 4463.6665
$ws.Range("K99").Value = 4463.6665
$ws.Range("M99").Value = -2965.6665

$ws.Range("H122").Value = 2800.7
$ws.Range("I122").Value = 2716.8572
$ws.Range("J122").Value = 2996.3333
$ws.Range("K122").Value = 8150.571599999999
$ws.Range("L122").Value = 8988.999899999999
$ws.Range("M122").Value = -5700.571599999999
$ws.Range("N122").Value = -13888.9999

$ws.Range("H126").Value = 4511.933
$ws.Range("I126").Value = 4463.6665
$ws.Range("K126").Value = 13390.9995
$ws.Range("M126").Value = -10920.9995

$ws.Range("H136").Value = 1412.625
$ws.Range("I136").Value = 786.4138
$ws.Range("K136").Value = 2359.2414
$ws.Range("M136").Value = 190.7586000000001

$ws.Range("H138").Value = 66373.60000000001
$ws.Range("J138").Value = 66373.60000000001
$ws.Range("L138").Value = 66373.60000000001
$ws.Range("N138").Value = -76653.60000000001

$ws.Range("H139").Value = 84990
$ws.Range("J139").Value = 84990
$ws.Range("L139").Value = 84990
$ws.Range("N139").Value = -95270

$ws.Range("H141").Value = 43393.3
$ws.Range("J141").Value = 43393.3
$ws.Range("L141").Value = 43393.3
$ws.Range("N141").Value = -53753.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 91.83871000000001
$ws.Range("I2").Value = 67.80952499999999
$ws.Range("J2").Value = 142.3
$ws.Range("K2").Value = 406.8571499999999
$ws.Range("L2").Value = 853.8000000000001
$ws.Range("M2").Value = -293.8571499999999
$ws.Range("N2").Value = -1079.8

$ws.Range("H36").Value = 533.6667
$ws.Range("I36").Value = 533.6667
$ws.Range("K36").Value = 1601.0001
$ws.Range("M36").Value = -1432.0001

$ws.Range("H129").Value = 2859228
$ws.Range("I129").Value = 6668087.5
$ws.Range("J129").Value = 2583.25
$ws.Range("K129").Value = 20004262.5
$ws.Range("L129").Value = 7749.75
$ws.Range("M129").Value = -19999262.5
$ws.Range("N129").Value = -17749.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10000
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws.Range("H99").Value = 2454
$ws.Range("I99").Value = 2454
$ws.Range("K99").Value = 2454
$ws.Range("M99").Value = -208

$ws.Range("H122").Value = 1115760.2
$ws.Range("I122").Value = 4446444
$ws.Range("J122").Value = 5532.3335
$ws.Range("K122").Value = 13339332
$ws.Range("L122").Value = 16597.0005
$ws.Range("M122").Value = -13336882
$ws.Range("N122").Value = -21497.0005

$ws.Range("H134").Value = 43752.777
$ws.Range("J134").Value = 43752.777
$ws.Range("L134").Value = 131258.331
$ws.Range("N134").Value = -136328.331

$ws.Range("H136").Value = 51802.816
$ws.Range("J136").Value = 51802.816
$ws.Range("L136").Value = 155408.448
$ws.Range("N136").Value = -160508.448

$ws.Range("H140").Value = 60166.668
$ws.Range("J140").Value = 60166.668
$ws.Range("L140").Value = 60166.668
$ws.Range("N140").Value = -70526.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6699999.5
$ws.Range("I2").Value = 10000000
$ws.Range("K2").Value = 10000000
$ws.Range("M2").Value = -9999888

$ws.Range("H122").Value = 5914.304
$ws.Range("I122").Value = 4006.5
$ws.Range("J122").Value = 7995.5454
$ws.Range("K122").Value = 12019.5
$ws.Range("L122").Value = 23986.6362
$ws.Range("M122").Value = -9569.5
$ws.Range("N122").Value = -28886.6362

$ws.Range("H132").Value = 4472.0386
$ws.Range("I132").Value = 4110.7954
$ws.Range("J132").Value = 6458.875
$ws.Range("K132").Value = 12332.3862
$ws.Range("L132").Value = 19376.625
$ws.Range("M132").Value = -9802.386200000001
$ws.Range("N132").Value = -24436.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H81").Value = 20836132
$ws.Range("I81").Value = 41669388
$ws.Range("J81").Value = 2875
$ws.Range("K81").Value = 83338776
$ws.Range("L81").Value = 5750
$ws.Range("M81").Value = -83337715
$ws.Range("N81").Value = -7872

$ws.Range("H84").Value = 20836132
$ws.Range("I84").Value = 41669388
$ws.Range("J84").Value = 2875
$ws.Range("K84").Value = 416693880
$ws.Range("L84").Value = 28750
$ws.Range("M84").Value = -416688576
$ws.Range("N84").Value = -39358

$ws.Range("H96").Value = 3214.5
$ws.Range("I96").Value = 3227.6365
$ws.Range("J96").Value = 3166.3333
$ws.Range("K96").Value = 3227.6365
$ws.Range("L96").Value = 3166.3333
$ws.Range("M96").Value = -1854.6365
$ws.Range("N96").Value = -5912.3333

$ws.Range("H122").Value = 3402.348
$ws.Range("I122").Value = 2692.1667
$ws.Range("J122").Value = 5959
$ws.Range("K122").Value = 8076.500100000001
$ws.Range("L122").Value = 17877
$ws.Range("M122").Value = -5626.500100000001
$ws.Range("N122").Value = -22777

$ws.Range("H137").Value = 98388.5
$ws.Range("J137").Value = 98388.5
$ws.Range("L137").Value = 98388.5
$ws.Range("N137").Value = -108588.5
